$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update stock-list cells (columns B-F) for rows 2-19 with the new ticker
# symbols. Only cells whose content actually changes are touched.
$ws.Range("B2").Value = "NSE:AMBUJACEM"
$ws.Range("C2").Value = "NSE:ALKALI"
$ws.Range("D2").Value = "NSE:BATAINDIA"
$ws.Range("E2").Value = "NSE:BEL"
$ws.Range("F2").Value = "NSE:BATAINDIA"
$ws.Range("B3").Value = "NSE:ANURAS"
$ws.Range("C3").Value = "NSE:ASHAPURMIN"
$ws.Range("D3").Value = "NSE:INFY"
$ws.Range("E3").Value = "NSE:EXIDEIND"
$ws.Range("F3").Value = "NSE:BHARTIARTL"
$ws.Range("B4").Value = "NSE:BHAGCHEM"
$ws.Range("C4").Value = "NSE:ASTEC"
$ws.Range("E4").Value = "NSE:HAL"
$ws.Range("F4").Value = "NSE:DLF"
$ws.Range("B5").Value = "NSE:COMPUSOFT"
$ws.Range("C5").Value = "NSE:BLAL"
$ws.Range("E5").Value = "NSE:IRCTC"
$ws.Range("F5").Value = "NSE:MARICO"
$ws.Range("B6").Value = "NSE:ESABINDIA"
$ws.Range("C6").Value = "NSE:BLUEJET"
$ws.Range("E6").Value = "NSE:LTTS"
$ws.Range("B7").Value = "NSE:ESG"
$ws.Range("C7").Value = "NSE:DMCC"
$ws.Range("E7").Value = "NSE:NMDC"
$ws.Range("B8").Value = "NSE:FINPIPE"
$ws.Range("C8").Value = "NSE:EXPLEOSOL"
$ws.Range("E8").Value = "NSE:NTPC"
$ws.Range("B9").Value = "NSE:INDUSTOWER"
$ws.Range("C9").Value = "NSE:GULPOLY"
$ws.Range("E9").Value = "NSE:PFC"
$ws.Range("B10").Value = "NSE:JKPAPER"
$ws.Range("C10").Value = "NSE:HAVISHA"
$ws.Range("E10").Value = "NSE:PVRINOX"
$ws.Range("B11").Value = "NSE:KEC"
$ws.Range("C11").Value = "NSE:HPL"
$ws.Range("E11").Value = "NSE:SAIL"
$ws.Range("B12").Value = "NSE:MANINFRA"
$ws.Range("C12").Value = "NSE:MADRASFERT"
$ws.Range("B13").Value = "NSE:MANOMAY"
$ws.Range("C13").Value = "NSE:MALLCOM"
$ws.Range("B14").Value = "NSE:MARICO"
$ws.Range("C14").Value = "NSE:MEGASTAR"
$ws.Range("B15").Value = "NSE:NAZARA"
$ws.Range("C15").Value = "NSE:NRBBEARING"
$ws.Range("C16").Value = "NSE:PANACEABIO"
$ws.Range("C17").Value = "NSE:POWERINDIA"
$ws.Range("C18").Value = "NSE:RACE"
$ws.Range("C19").Value = "NSE:RBL"

# The list shrank from 33 data rows to 18 data rows, so remove the
# now-unused trailing rows (this also shrinks the sheet dimension from
# A1:F34 down to A1:F19, matching the new data extent).
$ws.Rows("20:34").Delete()
